# Loop block shrinks and expands depending on height of child blocks
# Target slide is the last slide (slide 17) of the presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)

# ---------------------------------------------------------------------
# 1) Top "Group 17" (id 18) moves (no resize).
# ---------------------------------------------------------------------
$grp17 = $s.Shapes.Item("Group 17")
$grp17.Left = 402
$grp17.Top = 312

# ---------------------------------------------------------------------
# 2) "If" block header pieces move left (Rectangle 7, Rectangle 8,
#    Rectangle 10, Isosceles Triangle 11, Isosceles Triangle 12).
# ---------------------------------------------------------------------
$rect7 = $s.Shapes.Item("Rectangle 7")
$rect7.Left = 18

$rect8 = $s.Shapes.Item("Rectangle 8")
$rect8.Left = 56.75

$rect10 = $s.Shapes.Item("Rectangle 10")
$rect10.Left = 28.87496063

$tri11 = $s.Shapes.Item("Isosceles Triangle 11")
$tri11.Left = 18

$tri12 = $s.Shapes.Item("Isosceles Triangle 12")
$tri12.Left = 144

# ---------------------------------------------------------------------
# 3) Rectangle 13 / 14 / 15 / 16 become their own group ("Group 32"),
#    which is then repositioned as a whole.
# ---------------------------------------------------------------------
$loopRange = $s.Shapes.Range(@("Rectangle 13", "Rectangle 14", "Rectangle 15", "Rectangle 16"))
$loopGroup = $loopRange.Group()
$loopGroup.Name = "Group 32"
$loopGroup.Left = 18
$loopGroup.Top = 438

# ---------------------------------------------------------------------
# 4) "Group 23" moves.
# ---------------------------------------------------------------------
$grp23 = $s.Shapes.Item("Group 23")
$grp23.Left = 0

# ---------------------------------------------------------------------
# 5) "Group 33" moves.
# ---------------------------------------------------------------------
$grp33 = $s.Shapes.Item("Group 33")
$grp33.Left = 402
$grp33.Top = 126

# ---------------------------------------------------------------------
# 6) "Group 31" is ungrouped, its two rectangles are repositioned to
#    their new absolute (slide-space) locations and given explicit
#    fills (replacing the inherited grpFill), and then regrouped as
#    "Group 34". Doing the repositioning while the shapes are
#    ungrouped (top-level) makes the resulting group's transform an
#    identity translation (chOff == off), matching the target.
# ---------------------------------------------------------------------
$grp31 = $s.Shapes.Item("Group 31")
$grp31.Ungroup() | Out-Null

$rect28 = $s.Shapes.Item("Rectangle 28")
$rect29 = $s.Shapes.Item("Rectangle 29")

$rect28.Left = 402
$rect28.Top = 231
$rect28.Fill.ForeColor.RGB = 0x0033CC

$rect29.Left = 402 + 38.75
$rect29.Top = 231 + 60
$rect29.Fill.ForeColor.RGB = 0x0033CC

$newGroup31Range = $s.Shapes.Range(@("Rectangle 28", "Rectangle 29"))
$newGroup31 = $newGroup31Range.Group()
$newGroup31.Name = "Group 34"

# ---------------------------------------------------------------------
# 7) "Rectangle 30" moves up and shrinks in height.
# ---------------------------------------------------------------------
$rect30 = $s.Shapes.Item("Rectangle 30")
$rect30.Left = 402
$rect30.Top = 198
$rect30.Width = 36
$rect30.Height = 33
